$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 11000
$ws.Range("J54").Value = 11000
$ws.Range("L54").Value = 11000
$ws.Range("N54").Value = -11972
$ws.Range("H100").Value = 11906270
$ws.Range("I100").Value = 20834774
$ws.Range("J100").Value = 1600
$ws.Range("K100").Value = 20834774
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -20834233
$ws.Range("N100").Value = -2682
$ws.Range("H125").Value = 3515.0908
$ws.Range("I125").Value = 5997.75
$ws.Range("J125").Value = 2096.4285
$ws.Range("K125").Value = 53979.75
$ws.Range("L125").Value = 18867.8565
$ws.Range("M125").Value = -51519.75
$ws.Range("N125").Value = -23787.8565
$ws.Range("H131").Value = 2704.16
$ws.Range("I131").Value = 899
$ws.Range("J131").Value = 3155.45
$ws.Range("K131").Value = 2697
$ws.Range("L131").Value = 9466.349999999999
$ws.Range("M131").Value = 2343
$ws.Range("N131").Value = -19546.35
$ws.Range("H132").Value = 1712.5454
$ws.Range("I132").Value = 1103.9429
$ws.Range("J132").Value = 4079.3333
$ws.Range("K132").Value = 3311.8287
$ws.Range("L132").Value = 12237.9999
$ws.Range("M132").Value = -781.8287
$ws.Range("N132").Value = -17297.9999
$ws.Range("H137").Value = 2287.625
$ws.Range("I137").Value = 1935.8572
$ws.Range("J137").Value = 4750
$ws.Range("K137").Value = 5807.571599999999
$ws.Range("L137").Value = 14250
$ws.Range("M137").Value = -3257.571599999999
$ws.Range("N137").Value = -19350

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 5039.3335
$ws.Range("I30").Value = 3954
$ws.Range("J30").Value = 7210
$ws.Range("K30").Value = 3954
$ws.Range("L30").Value = 7210
$ws.Range("M30").Value = -3804
$ws.Range("N30").Value = -7510
$ws.Range("H35").Value = 4000
$ws.Range("I35").Value = 4000
$ws.Range("K35").Value = 4000
$ws.Range("M35").Value = -3594
$ws.Range("H102").Value = 3705115.5
$ws.Range("I102").Value = 4116593.8
$ws.Range("K102").Value = 4116593.8
$ws.Range("M102").Value = -4114971.8
$ws.Range("H118").Value = 32533.334
$ws.Range("J118").Value = 32533.334
$ws.Range("L118").Value = 32533.334
$ws.Range("N118").Value = -35847.334
$ws.Range("H123").Value = 29500
$ws.Range("J123").Value = 29500
$ws.Range("L123").Value = 29500
$ws.Range("N123").Value = -39300

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 9530
$ws.Range("J32").Value = 9530
$ws.Range("L32").Value = 9530
$ws.Range("N32").Value = -10298
$ws.Range("H38").Value = 11000
$ws.Range("J38").Value = 11000
$ws.Range("L38").Value = 11000
$ws.Range("N38").Value = -11832
$ws.Range("H99").Value = 55556964
$ws.Range("I99").Value = 125000744
$ws.Range("J99").Value = 1937.7
$ws.Range("K99").Value = 125000744
$ws.Range("L99").Value = 1937.7
$ws.Range("M99").Value = -124999246
$ws.Range("N99").Value = -4933.7
$ws.Range("H134").Value = 8724.765
$ws.Range("I134").Value = 16931.572
$ws.Range("J134").Value = 2980
$ws.Range("K134").Value = 50794.716
$ws.Range("L134").Value = 8940
$ws.Range("M134").Value = -48259.716
$ws.Range("N134").Value = -14010

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2318161.5
$ws.Range("I122").Value = 3969126.5
$ws.Range("J122").Value = 6810.6
$ws.Range("K122").Value = 11907379.5
$ws.Range("L122").Value = 20431.8
$ws.Range("M122").Value = -11904929.5
$ws.Range("N122").Value = -25331.8
$ws.Range("H132").Value = 5209.6
$ws.Range("I132").Value = 5900
$ws.Range("J132").Value = 5037
$ws.Range("K132").Value = 17700
$ws.Range("L132").Value = 15111
$ws.Range("M132").Value = -15170
$ws.Range("N132").Value = -20171
$ws.Range("H134").Value = 3500.182
$ws.Range("I134").Value = 3551.074
$ws.Range("J134").Value = 3271.1667
$ws.Range("K134").Value = 10653.222
$ws.Range("L134").Value = 9813.500100000001
$ws.Range("M134").Value = -8118.222
$ws.Range("N134").Value = -14883.5001

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2600.5
$ws.Range("I25").Value = 201
$ws.Range("K25").Value = 603
$ws.Range("M25").Value = -434
$ws.Range("H30").Value = 2600.5
$ws.Range("I30").Value = 201
$ws.Range("K30").Value = 603
$ws.Range("M30").Value = -501
$ws.Range("H31").Value = 2120
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 2525
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 7575
$ws.Range("M31").Value = -1212
$ws.Range("N31").Value = -8151
$ws.Range("H58").Value = 3404.8
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3404.8
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 10214.4
$ws.Range("M58").Value = $null
$ws.Range("N58").Value = -10470.4
$ws.Range("H107").Value = 372.25
$ws.Range("I107").Value = 165.11111
$ws.Range("J107").Value = 638.5714
$ws.Range("K107").Value = 495.33333
$ws.Range("L107").Value = 1915.7142
$ws.Range("M107").Value = 1424.66667
$ws.Range("N107").Value = -5755.7142
$ws.Range("H113").Value = 1622127.8
$ws.Range("I113").Value = 1724644.8
$ws.Range("J113").Value = 1250503.5
$ws.Range("K113").Value = 5173934.4
$ws.Range("L113").Value = 3751510.5
$ws.Range("M113").Value = -5171764.4
$ws.Range("N113").Value = -3755850.5
$ws.Range("H131").Value = 1640361.1
$ws.Range("I131").Value = 14286168
$ws.Range("J131").Value = 1089.8148
$ws.Range("K131").Value = 42858504
$ws.Range("L131").Value = 3269.4444
$ws.Range("M131").Value = -42853464
$ws.Range("N131").Value = -13349.4444

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = $null
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").Value = $null
$ws.Range("H122").Value = 1803131.4
$ws.Range("I122").Value = 3243096.5
$ws.Range("K122").Value = 9729289.5
$ws.Range("M122").Value = -9726839.5
$ws.Range("H126").Value = 4949
$ws.Range("I126").Value = 9499.385
$ws.Range("J126").Value = 2977.1667
$ws.Range("K126").Value = 28498.155
$ws.Range("L126").Value = 8931.500100000001
$ws.Range("M126").Value = -26028.155
$ws.Range("N126").Value = -13871.5001

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = $null
$ws.Range("H7").Value = 54558.58
$ws.Range("I7").Value = 64413
$ws.Range("J7").Value = 2001.6666
$ws.Range("K7").Value = 64413
$ws.Range("L7").Value = 2001.6666
$ws.Range("M7").Value = -64301
$ws.Range("N7").Value = -2225.6666
$ws.Range("H40").Value = 71431410
$ws.Range("I40").Value = 71431410
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 71431410
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -71431274
$ws.Range("N40").Value = $null
$ws.Range("H126").Value = 54558.58
$ws.Range("I126").Value = 64413
$ws.Range("J126").Value = 2001.6666
$ws.Range("K126").Value = 193239
$ws.Range("L126").Value = 6004.9998
$ws.Range("M126").Value = -190769
$ws.Range("N126").Value = -10944.9998

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 934.6667
$ws.Range("I96").Value = 958.6957
$ws.Range("J96").Value = 879.4
$ws.Range("K96").Value = 958.6957
$ws.Range("L96").Value = 879.4
$ws.Range("M96").Value = 414.3043
$ws.Range("N96").Value = -3625.4
$ws.Range("H113").Value = 939.5
$ws.Range("J113").Value = 1637.375
$ws.Range("L113").Value = 4912.125
$ws.Range("N113").Value = -9252.125
$ws.Range("H122").Value = 1301.2142
$ws.Range("I122").Value = 1035.6666
$ws.Range("J122").Value = 1500.375
$ws.Range("K122").Value = 3106.9998
$ws.Range("L122").Value = 4501.125
$ws.Range("M122").Value = -656.9998
$ws.Range("N122").Value = -9401.125
$ws.Range("H132").Value = 1972.069
$ws.Range("I132").Value = 1674.0416
$ws.Range("J132").Value = 3402.6
$ws.Range("K132").Value = 5022.1248
$ws.Range("L132").Value = 10207.8
$ws.Range("M132").Value = -2492.1248
$ws.Range("N132").Value = -15267.8
$ws.Range("H136").Value = 2315.1633
$ws.Range("I136").Value = 2460.1035
$ws.Range("J136").Value = 2105
$ws.Range("K136").Value = 7380.310500000001
$ws.Range("L136").Value = 6315
$ws.Range("M136").Value = -4830.310500000001
$ws.Range("N136").Value = -11415
